$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.008049845695496
$ws.Range("B1").Value = 2.107403755187988
$ws.Range("C1").Value = 5.63422679901123
$ws.Range("D1").Value = 0.8671468496322632
$ws.Range("E1").Value = 0.9318394660949707
